# Auto-generated edit script: applies the numeric corrections described in the commit diff.
# 200 cell-value updates + 5 cell deletions across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 9335.333000000001  # H32
$ws.Cells.Item(32, 10).Value = 7299.25  # J32
$ws.Cells.Item(32, 12).Value = 7299.25  # L32
$ws.Cells.Item(32, 14).Value = -7951.25  # N32
$ws.Cells.Item(43, 8).Value = 7010.8335  # H43
$ws.Cells.Item(43, 10).Value = 6433.1  # J43
$ws.Cells.Item(43, 12).Value = 6433.1  # L43
$ws.Cells.Item(43, 14).Value = -6571.1  # N43
$ws.Cells.Item(76, 8).Value = 9000.4  # H76
$ws.Cells.Item(76, 9).Value = 4000  # I76
$ws.Cells.Item(76, 10).Value = 11143.429  # J76
$ws.Cells.Item(76, 11).Value = 4000  # K76
$ws.Cells.Item(76, 12).Value = 11143.429  # L76
$ws.Cells.Item(76, 13).Value = -3685  # M76
$ws.Cells.Item(76, 14).Value = -11773.429  # N76
$ws.Cells.Item(79, 8).Value = 9000.4  # H79
$ws.Cells.Item(79, 9).Value = 4000  # I79
$ws.Cells.Item(79, 10).Value = 11143.429  # J79
$ws.Cells.Item(79, 11).Value = 4000  # K79
$ws.Cells.Item(79, 12).Value = 11143.429  # L79
$ws.Cells.Item(79, 13).Value = -2908  # M79
$ws.Cells.Item(79, 14).Value = -13327.429  # N79
$ws.Cells.Item(98, 8).Value = 8367.166999999999  # H98
$ws.Cells.Item(98, 9).Value = 534  # I98
$ws.Cells.Item(98, 10).Value = 12283.75  # J98
$ws.Cells.Item(98, 11).Value = 534  # K98
$ws.Cells.Item(98, 12).Value = 12283.75  # L98
$ws.Cells.Item(98, 13).Value = 964  # M98
$ws.Cells.Item(98, 14).Value = -15279.75  # N98
$ws.Cells.Item(106, 8).Value = 12841.556  # H106
$ws.Cells.Item(106, 9).Value = 5471.1665  # I106
$ws.Cells.Item(106, 10).Value = 16526.75  # J106
$ws.Cells.Item(106, 11).Value = 5471.1665  # K106
$ws.Cells.Item(106, 12).Value = 16526.75  # L106
$ws.Cells.Item(106, 13).Value = -4840.1665  # M106
$ws.Cells.Item(106, 14).Value = -17788.75  # N106
$ws.Cells.Item(122, 8).Value = 8367.166999999999  # H122
$ws.Cells.Item(122, 9).Value = 534  # I122
$ws.Cells.Item(122, 10).Value = 12283.75  # J122
$ws.Cells.Item(122, 11).Value = 1602  # K122
$ws.Cells.Item(122, 12).Value = 36851.25  # L122
$ws.Cells.Item(122, 13).Value = 848  # M122
$ws.Cells.Item(122, 14).Value = -41751.25  # N122
$ws.Cells.Item(138, 8).Value = 3836.5264  # H138
$ws.Cells.Item(138, 10).Value = 3676.923  # J138
$ws.Cells.Item(138, 12).Value = 11030.769  # L138
$ws.Cells.Item(138, 14).Value = -21310.769  # N138

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4179.5386  # H45
$ws.Cells.Item(45, 9).Value = 1742  # I45
$ws.Cells.Item(45, 11).Value = 1742  # K45
$ws.Cells.Item(45, 13).Value = -1365  # M45
$ws.Cells.Item(74, 8).Value = 33335808  # H74
$ws.Cells.Item(74, 9).Value = 33335808  # I74
$ws.Cells.Item(74, 10).Value = 0  # J74
$ws.Cells.Item(74, 11).Value = 33335808  # K74
$ws.Cells.Item(74, 12).Value = 0  # L74
$ws.Cells.Item(74, 13).Value = -33334934  # M74
$ws.Cells.Item(77, 8).Value = 33335808  # H77
$ws.Cells.Item(77, 9).Value = 33335808  # I77
$ws.Cells.Item(77, 10).Value = 0  # J77
$ws.Cells.Item(77, 11).Value = 166679040  # K77
$ws.Cells.Item(77, 12).Value = 0  # L77
$ws.Cells.Item(77, 13).Value = -166674672  # M77
$ws.Cells.Item(108, 8).Value = 0  # H108
$ws.Cells.Item(108, 10).Value = 0  # J108
$ws.Cells.Item(108, 12).Value = 0  # L108
$ws.Cells.Item(111, 8).Value = 65000  # H111
$ws.Cells.Item(111, 10).Value = 65000  # J111
$ws.Cells.Item(111, 12).Value = 65000  # L111
$ws.Cells.Item(111, 14).Value = -73180  # N111
$ws.Cells.Item(113, 8).Value = 63949  # H113
$ws.Cells.Item(113, 10).Value = 63949  # J113
$ws.Cells.Item(113, 12).Value = 63949  # L113
$ws.Cells.Item(113, 14).Value = -72627  # N113
$ws.Cells.Item(119, 8).Value = 71841.5  # H119
$ws.Cells.Item(119, 10).Value = 71841.5  # J119
$ws.Cells.Item(119, 12).Value = 71841.5  # L119
$ws.Cells.Item(119, 14).Value = -81517.5  # N119
$ws.Cells.Item(132, 8).Value = 2218.4531  # H132
$ws.Cells.Item(132, 9).Value = 1515.4386  # I132
$ws.Cells.Item(132, 11).Value = 4546.3158  # K132
$ws.Cells.Item(132, 13).Value = -2016.3158  # M132
$ws.Cells.Item(74, 14).ClearContents()  # N74 (cell removed)
$ws.Cells.Item(77, 14).ClearContents()  # N77 (cell removed)
$ws.Cells.Item(108, 14).ClearContents()  # N108 (cell removed)

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 1792.4445  # H80
$ws.Cells.Item(80, 10).Value = 2191.5  # J80
$ws.Cells.Item(80, 12).Value = 2191.5  # L80
$ws.Cells.Item(80, 14).Value = -4187.5  # N80
$ws.Cells.Item(83, 8).Value = 1792.4445  # H83
$ws.Cells.Item(83, 10).Value = 2191.5  # J83
$ws.Cells.Item(83, 12).Value = 10957.5  # L83
$ws.Cells.Item(83, 14).Value = -20941.5  # N83
$ws.Cells.Item(105, 8).Value = 37763.375  # H105
$ws.Cells.Item(105, 9).Value = 52402.25  # I105
$ws.Cells.Item(105, 10).Value = 23124.5  # J105
$ws.Cells.Item(105, 11).Value = 52402.25  # K105
$ws.Cells.Item(105, 12).Value = 23124.5  # L105
$ws.Cells.Item(105, 13).Value = -50655.25  # M105
$ws.Cells.Item(105, 14).Value = -26618.5  # N105
$ws.Cells.Item(108, 8).Value = 0  # H108
$ws.Cells.Item(108, 10).Value = 0  # J108
$ws.Cells.Item(108, 12).Value = 0  # L108
$ws.Cells.Item(108, 14).ClearContents()  # N108 (cell removed)

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 150001  # H6
$ws.Cells.Item(6, 9).Value = 150001  # I6
$ws.Cells.Item(6, 11).Value = 150001  # K6
$ws.Cells.Item(6, 13).Value = -149888  # M6
$ws.Cells.Item(31, 8).Value = 57115.383  # H31
$ws.Cells.Item(31, 9).Value = 7667.5  # I31
$ws.Cells.Item(31, 11).Value = 7667.5  # K31
$ws.Cells.Item(31, 13).Value = -7372.5  # M31
$ws.Cells.Item(34, 8).Value = 57115.383  # H34
$ws.Cells.Item(34, 9).Value = 7667.5  # I34
$ws.Cells.Item(34, 11).Value = 7667.5  # K34
$ws.Cells.Item(34, 13).Value = -7465.5  # M34
$ws.Cells.Item(94, 8).Value = 3555.75  # H94
$ws.Cells.Item(94, 9).Value = 1733.25  # I94
$ws.Cells.Item(94, 11).Value = 1733.25  # K94
$ws.Cells.Item(94, 13).Value = -1282.25  # M94
$ws.Cells.Item(119, 8).Value = 0  # H119
$ws.Cells.Item(119, 10).Value = 0  # J119
$ws.Cells.Item(119, 12).Value = 0  # L119
$ws.Cells.Item(121, 8).Value = 41249.75  # H121
$ws.Cells.Item(121, 10).Value = 41249.75  # J121
$ws.Cells.Item(121, 12).Value = 41249.75  # L121
$ws.Cells.Item(121, 14).Value = -43869.75  # N121
$ws.Cells.Item(134, 8).Value = 3352.8823  # H134
$ws.Cells.Item(134, 9).Value = 1891.2307  # I134
$ws.Cells.Item(134, 10).Value = 8103.25  # J134
$ws.Cells.Item(134, 11).Value = 5673.6921  # K134
$ws.Cells.Item(134, 12).Value = 24309.75  # L134
$ws.Cells.Item(134, 13).Value = -3138.6921  # M134
$ws.Cells.Item(134, 14).Value = -29379.75  # N134
$ws.Cells.Item(119, 14).ClearContents()  # N119 (cell removed)

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1667.3077  # H107
$ws.Cells.Item(107, 10).Value = 1645.4667  # J107
$ws.Cells.Item(107, 12).Value = 4936.4001  # L107
$ws.Cells.Item(107, 14).Value = -8776.400099999999  # N107

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 9830.5  # H80
$ws.Cells.Item(80, 9).Value = 8080  # I80
$ws.Cells.Item(80, 10).Value = 11581  # J80
$ws.Cells.Item(80, 11).Value = 8080  # K80
$ws.Cells.Item(80, 12).Value = 11581  # L80
$ws.Cells.Item(80, 13).Value = -7082  # M80
$ws.Cells.Item(80, 14).Value = -13577  # N80
$ws.Cells.Item(83, 8).Value = 9830.5  # H83
$ws.Cells.Item(83, 9).Value = 8080  # I83
$ws.Cells.Item(83, 10).Value = 11581  # J83
$ws.Cells.Item(83, 11).Value = 40400  # K83
$ws.Cells.Item(83, 12).Value = 57905  # L83
$ws.Cells.Item(83, 13).Value = -35408  # M83
$ws.Cells.Item(83, 14).Value = -67889  # N83
$ws.Cells.Item(102, 8).Value = 3038.318  # H102
$ws.Cells.Item(102, 9).Value = 2122.8667  # I102
$ws.Cells.Item(102, 11).Value = 2122.8667  # K102
$ws.Cells.Item(102, 13).Value = -500.8667  # M102
$ws.Cells.Item(107, 8).Value = 6000  # H107
$ws.Cells.Item(107, 10).Value = 6000  # J107
$ws.Cells.Item(107, 12).Value = 6000  # L107
$ws.Cells.Item(107, 14).Value = -9840  # N107
$ws.Cells.Item(126, 8).Value = 3347.7368  # H126
$ws.Cells.Item(126, 9).Value = 2477.389  # I126
$ws.Cells.Item(126, 11).Value = 7432.167  # K126
$ws.Cells.Item(126, 13).Value = -4962.167  # M126
$ws.Cells.Item(132, 8).Value = 4650  # H132
$ws.Cells.Item(132, 9).Value = 3816.2144  # I132
$ws.Cells.Item(132, 10).Value = 10486.5  # J132
$ws.Cells.Item(132, 11).Value = 11448.6432  # K132
$ws.Cells.Item(132, 12).Value = 31459.5  # L132
$ws.Cells.Item(132, 13).Value = -8918.643199999999  # M132
$ws.Cells.Item(132, 14).Value = -36519.5  # N132

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4896.8076  # H22
$ws.Cells.Item(22, 9).Value = 1670.7273  # I22
$ws.Cells.Item(22, 10).Value = 7262.6  # J22
$ws.Cells.Item(22, 11).Value = 1670.7273  # K22
$ws.Cells.Item(22, 12).Value = 7262.6  # L22
$ws.Cells.Item(22, 13).Value = -1375.7273  # M22
$ws.Cells.Item(22, 14).Value = -7852.6  # N22
$ws.Cells.Item(27, 8).Value = 4896.8076  # H27
$ws.Cells.Item(27, 9).Value = 1670.7273  # I27
$ws.Cells.Item(27, 10).Value = 7262.6  # J27
$ws.Cells.Item(27, 11).Value = 1670.7273  # K27
$ws.Cells.Item(27, 12).Value = 7262.6  # L27
$ws.Cells.Item(27, 13).Value = -1563.7273  # M27
$ws.Cells.Item(27, 14).Value = -7476.6  # N27
$ws.Cells.Item(45, 8).Value = 13500  # H45
$ws.Cells.Item(45, 10).Value = 13500  # J45
$ws.Cells.Item(45, 12).Value = 13500  # L45
$ws.Cells.Item(45, 14).Value = -14314  # N45
$ws.Cells.Item(46, 8).Value = 2616.1667  # H46
$ws.Cells.Item(46, 10).Value = 2959.6  # J46
$ws.Cells.Item(46, 12).Value = 2959.6  # L46
$ws.Cells.Item(46, 14).Value = -3335.6  # N46
$ws.Cells.Item(136, 8).Value = 7473.75  # H136
$ws.Cells.Item(136, 9).Value = 1674.2727  # I136
$ws.Cells.Item(136, 10).Value = 14562  # J136
$ws.Cells.Item(136, 11).Value = 5022.8181  # K136
$ws.Cells.Item(136, 12).Value = 43686  # L136
$ws.Cells.Item(136, 13).Value = -2472.8181  # M136
$ws.Cells.Item(136, 14).Value = -48786  # N136

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 6425.34  # H132
$ws.Cells.Item(132, 9).Value = 4427.3335  # I132
$ws.Cells.Item(132, 11).Value = 13282.0005  # K132
$ws.Cells.Item(132, 13).Value = -10752.0005  # M132
